# Replace the CJK placeholder text (and slide-layout display names) in the
# slide master / slide layouts with plain ASCII strings, to get consistent
# fonts used - mirrors the upstream LibreOffice commit for
# sd/qa/unit/data/pptx/smartart-tdf134221.pptx.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

function Set-ParaTexts($shape, $texts) {
    $tr = $shape.TextFrame.TextRange
    for ($i = 0; $i -lt $texts.Count; $i++) {
        $tr.Paragraphs($i + 1, 1).Text = $texts[$i]
    }
}

# --- slideLayout1.xml (type="title") -------------------------------------
$cl = $master.CustomLayouts.Item(1)
$cl.Name = "THING"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(2).TextFrame.TextRange.Text = "WORLD"

# --- slideLayout2.xml (type="obj") ----------------------------------------
$cl = $master.CustomLayouts.Item(2)
$cl.Name = "STYLE"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-ParaTexts $cl.Shapes.Item(2) @("SOMETHING","FOO","BAR","BAR","BAZ")

# --- slideLayout3.xml (type="secHead") ------------------------------------
$cl = $master.CustomLayouts.Item(3)
$cl.Name = "ONE"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(2).TextFrame.TextRange.Text = "SOMETHING"

# --- slideLayout4.xml (type="twoObj") -------------------------------------
$cl = $master.CustomLayouts.Item(4)
$cl.Name = "FOUR"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-ParaTexts $cl.Shapes.Item(2) @("SOMETHING","FOO","BAR","BAR","BAZ")
Set-ParaTexts $cl.Shapes.Item(3) @("SOMETHING","FOO","BAR","BAR","BAZ")

# --- slideLayout5.xml (type="twoTxTwoObj") --------------------------------
$cl = $master.CustomLayouts.Item(5)
$cl.Name = "EG"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(2).TextFrame.TextRange.Text = "SOMETHING"
Set-ParaTexts $cl.Shapes.Item(3) @("SOMETHING","FOO","BAR","BAR","BAZ")
$cl.Shapes.Item(4).TextFrame.TextRange.Text = "SOMETHING"
Set-ParaTexts $cl.Shapes.Item(5) @("SOMETHING","FOO","BAR","BAR","BAZ")

# --- slideLayout6.xml (type="titleOnly") ----------------------------------
$cl = $master.CustomLayouts.Item(6)
$cl.Name = "ABC"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"

# --- slideLayout7.xml (type="blank") ---------------------------------------
$cl = $master.CustomLayouts.Item(7)
$cl.Name = "XY"

# --- slideLayout8.xml (type="objTx") ---------------------------------------
$cl = $master.CustomLayouts.Item(8)
$cl.Name = "ABCDEF"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-ParaTexts $cl.Shapes.Item(2) @("SOMETHING","FOO","BAR","BAR","BAZ")
$cl.Shapes.Item(3).TextFrame.TextRange.Text = "SOMETHING"

# --- slideLayout9.xml (type="picTx") ----------------------------------------
$cl = $master.CustomLayouts.Item(9)
$cl.Name = "HIJKL"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(3).TextFrame.TextRange.Text = "SOMETHING"

# --- slideLayout10.xml (type="vertTx") --------------------------------------
$cl = $master.CustomLayouts.Item(10)
$cl.Name = "EXAMPLE"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-ParaTexts $cl.Shapes.Item(2) @("SOMETHING","FOO","BAR","BAR","BAZ")

# --- slideLayout11.xml (type="vertTitleAndTx") ------------------------------
$cl = $master.CustomLayouts.Item(11)
$cl.Name = "SOMETHING"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-ParaTexts $cl.Shapes.Item(2) @("SOMETHING","FOO","BAR","BAR","BAZ")

# --- slideMaster1.xml --------------------------------------------------------
$master.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-ParaTexts $master.Shapes.Item(2) @("SOMETHING","FOO","BAR","BAR","BAZ")
